# 20190503_plate_layout.xlsx edit
#
# Re-labels the well-coordinate strings on the "well" sheet to use
# zero-padded column numbers (e.g. "B2" -> "B02", "G9" -> "G09"); the
# already two-digit labels ("B10", "B11", ...) are left as-is. Also
# switches the active sheet/selection from "misc" to "well" (exploring
# the well grid next while using holoviews), and leaves the "misc"
# sheet's own remembered selection untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("well")

$rowLetters = @("B", "C", "D", "E", "F", "G")
$colLetters = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")

for ($r = 0; $r -lt $rowLetters.Length; $r++) {
    $letter = $rowLetters[$r]
    for ($c = 0; $c -lt $colLetters.Length; $c++) {
        $n = $c + 2
        if ($n -le 9) {
            $label = "{0}0{1}" -f $letter, $n
        } else {
            $label = "{0}{1}" -f $letter, $n
        }
        $cellRef = "{0}{1}" -f $colLetters[$c], ($r + 1)
        $ws.Range($cellRef).Value = $label
    }
}

# Make "well" the active sheet/tab, with I11 selected on it (this also
# clears the previous tabSelected flag + selection memory on "misc").
$ws.Activate()
$ws.Range("I11").Select() | Out-Null
